$d = $word.ActiveDocument

$replacements = @(
    @{old = "56×54=3024"; new = "37×51=1887"},
    @{old = "52×16=832"; new = "17×44=748"},
    @{old = "15×77=1155"; new = "25×13=325"},
    @{old = "69×42=2898"; new = "35×21=735"},
    @{old = "96×42=4032"; new = "71×59=4189"},
    @{old = "12×23=276"; new = "89×49=4361"},
    @{old = "31×20=620"; new = "92×24=2208"},
    @{old = "61×81=4941"; new = "19×45=855"},
    @{old = "25×49=1225"; new = "42×86=3612"},
    @{old = "16×82=1312"; new = "19×53=1007"},
    @{old = "33×30=990"; new = "42×88=3696"},
    @{old = "11×51=561"; new = "65×71=4615"},
    @{old = "47×73=3431"; new = "15×75=1125"},
    @{old = "57×55=3135"; new = "66×22=1452"},
    @{old = "28×56=1568"; new = "83×38=3154"},
    @{old = "60×12=720"; new = "68×74=5032"},
    @{old = "59×12=708"; new = "65×96=6240"},
    @{old = "42×41=1722"; new = "36×64=2304"},
    @{old = "83×20=1660"; new = "72×28=2016"},
    @{old = "70×70=4900"; new = "86×11=946"},
    @{old = "84×11=924"; new = "76×76=5776"},
    @{old = "46×93=4278"; new = "17×31=527"},
    @{old = "55×62=3410"; new = "28×51=1428"},
    @{old = "98×86=8428"; new = "88×57=5016"},
    @{old = "96×90=8640"; new = "16×73=1168"}
)

foreach ($pair in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
